# Add daily power records: extend the "comforter_cda_table" table from
# A1:F50 to A1:F51, backfill missing Start/End Time values for several
# existing rows, add a brand-new row 50 (date 43374) and row 51
# (date 43375), and let the Duration/Second Duration/Absolute Value
# formulas recompute (and re-group into a shared-formula block for
# D42:F50, matching how Excel re-shares formulas after a fill).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the Excel table by one row so its ref becomes A1:F51 (also moves
# the autoFilter range and the sheet dimension along with it).
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# --- Backfill Start Time / End Time for existing rows that had none ---
$ws.Range("B39").Value = 0
$ws.Range("C39").Value = 0.36319444444444443

$ws.Range("B40").Value = 0
$ws.Range("C40").Value = 0

$ws.Range("B43").Value = 0
$ws.Range("C43").Value = 0

$ws.Range("B44").Value = 0.43055555555555558
$ws.Range("C44").Value = 0.65069444444444446

$ws.Range("B45").Value = 0.75694444444444453
$ws.Range("C45").Value = 0.99930555555555556

$ws.Range("B46").Value = 0
$ws.Range("C46").Value = 0.40138888888888885

# Row 47 already had placeholder zeros - overwrite with real times.
$ws.Range("B47").Value = 0.46736111111111112
$ws.Range("C47").Value = 0.67986111111111114

# --- New row 50 data (date + times) ---
$ws.Range("A50").Value = 43374
$ws.Range("B50").Value = 0.81805555555555554
$ws.Range("C50").Value = 0.99930555555555556

# --- New row 51 - date only, Start/End Time left blank ---
$ws.Range("A51").Value = 43375

# Re-enter the calculated-column formulas across D42:F50 so Excel
# regroups them into a shared-formula block (matches rows 2-41 which
# already use t="shared" groups).
$ws.Range("D42:D50").Formula = "=(C42-B42)* 1440"
$ws.Range("E42:E50").Formula = "=IF(C42>B42, (C42-B42)*1440, (B42-C42)*1440)"
$ws.Range("F42:F50").Formula = "=ABS((C42-B42)*1440)"

# Row 51's formulas (new table row) stay as plain, non-shared formulas.
$ws.Range("D51").Formula = "=(C51-B51)* 1440"
$ws.Range("E51").Formula = "=IF(C51>B51, (C51-B51)*1440, (B51-C51)*1440)"
$ws.Range("F51").Formula = "=ABS((C51-B51)*1440)"

# Match the saved selection/scroll state from the authored edit.
$ws.Range("B51").Select() | Out-Null
